$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# The sheet already contains a "week" block in rows 211-220:
#   211       : merged header "日期：2018.11.26 第十三周周一"
#   212       : column titles (组员 / 计划内容 / 完成情况 / 备注)
#   213-218   : one row per team member
#   219-220   : merged "总结：..." footer
# A new, almost identical block needs to be appended as rows 221-230 for
# "日期：2018.11.28 第十三周周三" with every member's plan set to "项目整合"
# and an (still empty, to be filled in later) summary footer.

# 1) Duplicate the block's values/formulas/merges via a normal paste ...
$ws.Range("A211:D220").Copy()
$ws.Range("A221:D230").PasteSpecial()

# 2) ... then re-stamp the exact per-cell formatting on top (a plain paste
#    alone resets everything to the default style).
$ws.Range("A211:D220").Copy()
$ws.Range("A221:D230").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Overwrite the header text for the new block. The date part stays in
#    the default (regular) run while " 第十三周周三" is bold, matching the
#    formatting already used for the other weekly headers.
$headerDate = "日期：2018.11.28"
$headerWeek = " 第十三周周三"
$ws.Range("A221").Value = $headerDate + $headerWeek
$boldRun = $ws.Range("A221").Characters($headerDate.Length + 1, $headerWeek.Length)
$boldRun.Font.Bold = $true
$boldRun.Font.Size = 10

# 4) Every team member's plan content for this entry is "项目整合".
$ws.Range("B223").Value = "项目整合"
$ws.Range("B224").Value = "项目整合"
$ws.Range("B225").Value = "项目整合"
$ws.Range("B226").Value = "项目整合"
$ws.Range("B227").Value = "项目整合"
$ws.Range("B228").Value = "项目整合"

# Row 213 in the source block had leftover "完成情况/备注" text (C213/D213) -
# the new block's corresponding row (223) should start blank like the rest.
$ws.Range("C223").ClearContents()
$ws.Range("D223").ClearContents()

# 5) Match the sheet's reported selection state after the edit.
$ws.Range("A219:D220").Select()
